$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price values that look numeric (e.g. "572.96",
# "0.0000107") but must stay as literal text, matching the inline-string
# cells in the workbook. Force a Text number format before assigning the
# string, then reset the style back to Normal so no stray formatting is
# left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.499.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.597.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.68%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -4.37%  '
$ws.Range('E9').Value = '  -7.90%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.380'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.065.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('E15').Value = '  -9.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.353.65'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.585.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('E18').Value = '  -5.55%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '340.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.73%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000107'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '577.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E28').Value = '  -3.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('E32').Value = '  -4.82%  '
$ws.Range('E33').Value = '  -5.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.23%  '
$ws.Range('E36').Value = '  -5.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.69'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +6.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '156.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0587'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.627'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('E50').Value = '  -5.14%  '
$ws.Range('E51').Value = '  -5.59%  '
